$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.640.59"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.028.84"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.31"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.41"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.86"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "3.505.91"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.58"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.78"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "3.058.46"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.995"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.69"
$ws.Range("E18").Value = "  -10.76%  "
$ws.Range("D19").Value = "51.728.24"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.07"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.31"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.86"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.18"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.38"
$ws.Range("E26").Value = "  +7.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.48"
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.23"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.34"
$ws.Range("E32").Value = "  +4.42%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.23"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.34"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +8.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.30"
$ws.Range("E39").Value = "  +6.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.62"
$ws.Range("E40").Value = "  +9.24%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.282"
$ws.Range("E42").Value = "  +9.84%  "
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.94"
$ws.Range("E44").Value = "  +6.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.78"
$ws.Range("E45").Value = "  +13.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.86"
$ws.Range("E46").Value = "  +5.37%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").Value = "2.041.94"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("D50").Value = "3.337.36"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0323"
$ws.Range("E51").Value = "  +3.34%  "
